$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price value is a plain number-looking string need to be
# forced to Text format first, or Excel's COM Value setter will coerce them
# into numeric cells (the source data models price as text, e.g. '404.55').
$forceTextCells = @("D5","D6","D7","D9","D11","D14","D15","D17","D21","D22","D23","D24","D25","D26","D27","D28","D29","D34","D37","D38","D40","D41","D43","D44","D46","D48","D51")
foreach ($ref in $forceTextCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = '61.127.13'
$ws.Range("E2").Value = '  -0.09%  '

$ws.Range("D3").Value = '3.382.05'
$ws.Range("E3").Value = '  +1.98%  '

$ws.Range("E4").Value = '  +0.15%  '

$ws.Range("D5").Value = '404.55'
$ws.Range("E5").Value = '  -1.43%  '

$ws.Range("D6").Value = '128.38'
$ws.Range("E6").Value = '  +14.02%  '

$ws.Range("D7").Value = '0.609'
$ws.Range("E7").Value = '  +7.44%  '

$ws.Range("E8").Value = '  +0.16%  '

$ws.Range("D9").Value = '0.675'
$ws.Range("E9").Value = '  +8.71%  '

$ws.Range("E10").Value = '  +9.71%  '

$ws.Range("D11").Value = '42.16'
$ws.Range("E11").Value = '  +8.90%  '

$ws.Range("E12").Value = '  -0.47%  '

$ws.Range("D13").Value = '3.931.32'
$ws.Range("E13").Value = '  +2.32%  '

$ws.Range("D14").Value = '8.51'
$ws.Range("E14").Value = '  +4.47%  '

$ws.Range("D15").Value = '19.65'
$ws.Range("E15").Value = '  +3.41%  '

$ws.Range("D16").Value = '3.379.23'
$ws.Range("E16").Value = '  +1.76%  '

$ws.Range("D17").Value = '11.47'
$ws.Range("E17").Value = '  +8.77%  '

$ws.Range("D18").Value = '61.085.75'
$ws.Range("E18").Value = '  +0.25%  '

$ws.Range("E19").Value = '  +3.77%  '

$ws.Range("E20").Value = '  +17.79%  '

$ws.Range("D21").Value = '3.25'
$ws.Range("E21").Value = '  +1.01%  '

$ws.Range("D22").Value = '82.79'
$ws.Range("E22").Value = '  +13.21%  '

$ws.Range("D23").Value = '13.13'
$ws.Range("E23").Value = '  +6.75%  '

$ws.Range("D24").Value = '306.66'
$ws.Range("E24").Value = '  +4.20%  '

$ws.Range("D25").Value = '3.15'
$ws.Range("E25").Value = '  +2.36%  '

$ws.Range("B26").Value = 'Filecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D26").Value = '8.50'
$ws.Range("E26").Value = '  +14.72%  '

$ws.Range("B27").Value = 'LEO'
$ws.Range("C27").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D27").Value = '4.72'
$ws.Range("E27").Value = '  +4.50%  '

$ws.Range("D28").Value = '29.60'
$ws.Range("E28").Value = '  +2.78%  '

$ws.Range("D29").Value = '7.46'
$ws.Range("E29").Value = '  +1.70%  '

$ws.Range("E30").Value = '  +0.53%  '

$ws.Range("E31").Value = '  +6.63%  '

$ws.Range("E32").Value = '  +6.38%  '

$ws.Range("E33").Value = '  +5.98%  '

$ws.Range("D34").Value = '42.39'
$ws.Range("E34").Value = '  +6.81%  '

$ws.Range("E35").Value = '  +0.17%  '

$ws.Range("E36").Value = '  +1.52%  '

$ws.Range("D37").Value = '52.28'
$ws.Range("E37").Value = '  -0.46%  '

$ws.Range("D38").Value = '0.997'
$ws.Range("E38").Value = '  -0.12%  '

$ws.Range("E39").Value = '  +4.55%  '

$ws.Range("D40").Value = '2.98'
$ws.Range("E40").Value = '  -1.65%  '

$ws.Range("D41").Value = '2.02'
$ws.Range("E41").Value = '  +7.96%  '

$ws.Range("E42").Value = '  +5.23%  '

$ws.Range("D43").Value = '137.06'
$ws.Range("E43").Value = '  +1.25%  '

$ws.Range("D44").Value = '3.96'
$ws.Range("E44").Value = '  +5.63%  '

$ws.Range("E45").Value = '  +0.90%  '

$ws.Range("D46").Value = '16.94'
$ws.Range("E46").Value = '  +4.85%  '

$ws.Range("E47").Value = '  +1.39%  '

$ws.Range("D48").Value = '21.73'
$ws.Range("E48").Value = '  +4.29%  '

$ws.Range("D49").Value = '2.140.84'
$ws.Range("E49").Value = '  +1.42%  '

$ws.Range("D50").Value = '3.715.40'
$ws.Range("E50").Value = '  +1.89%  '

$ws.Range("D51").Value = '2.35'
$ws.Range("E51").Value = '  +0.91%  '

# Restore the default (Normal) style on the forced-text cells so no stray
# number-format style lingers on them once the text value is committed.
foreach ($ref in $forceTextCells) {
    $ws.Range($ref).Style = "Normal"
}
